{"js": "// Update the division-problem answers in the results table.\n// Each entry is a (old answer text -> new answer text) pair; every old\n// string occurs exactly once in the document, so a body.search() +\n// insertText(replace) safely retargets only the intended cell.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"999\u00f74=249, 3\", \"269\u00f75=53, 4\"],\n  [\"829\u00f75=165, 4\", \"748\u00f74=187, 0\"],\n  [\"793\u00f79=88, 1\", \"823\u00f74=205, 3\"],\n  [\"620\u00f77=88, 4\", \"985\u00f77=140, 5\"],\n  [\"249\u00f76=41, 3\", \"882\u00f75=176, 2\"],\n  [\"132\u00f75=26, 2\", \"624\u00f72=312, 0\"],\n  [\"623\u00f78=77, 7\", \"595\u00f72=297, 1\"],\n  [\"528\u00f79=58, 6\", \"779\u00f73=259, 2\"],\n  [\"465\u00f75=93, 0\", \"889\u00f77=127, 0\"],\n  [\"335\u00f78=41, 7\", \"979\u00f73=326, 1\"],\n  [\"313\u00f76=52, 1\", \"855\u00f79=95, 0\"],\n  [\"956\u00f77=136, 4\", \"104\u00f77=14, 6\"],\n  [\"442\u00f74=110, 2\", \"388\u00f72=194, 0\"],\n  [\"227\u00f79=25, 2\", \"391\u00f72=195, 1\"],\n  [\"471\u00f76=78, 3\", \"395\u00f79=43, 8\"],\n  [\"267\u00f79=29, 6\", \"946\u00f72=473, 0\"],\n  [\"858\u00f79=95, 3\", \"825\u00f72=412, 1\"],\n  [\"583\u00f73=194, 1\", \"690\u00f79=76, 6\"],\n  [\"564\u00f72=282, 0\", \"995\u00f77=142, 1\"],\n  [\"824\u00f76=137, 2\", \"855\u00f79=95, 0\"],\n  [\"390\u00f76=65, 0\", \"503\u00f76=83, 5\"],\n  [\"290\u00f79=32, 2\", \"202\u00f72=101, 0\"],\n  [\"580\u00f72=290, 0\", \"546\u00f78=68, 2\"],\n  [\"556\u00f73=185, 1\", \"329\u00f79=36, 5\"],\n  [\"303\u00f74=75, 3\", \"591\u00f78=73, 7\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem answers in the results table.\n# Each entry is a (old answer text -> new answer text) pair; every old\n# string occurs exactly once in the document, so a plain Find/Replace\n# (no wildcards) safely retargets only the intended cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"999\u00f74=249, 3\", \"269\u00f75=53, 4\"),\n    @(\"829\u00f75=165, 4\", \"748\u00f74=187, 0\"),\n    @(\"793\u00f79=88, 1\", \"823\u00f74=205, 3\"),\n    @(\"620\u00f77=88, 4\", \"985\u00f77=140, 5\"),\n    @(\"249\u00f76=41, 3\", \"882\u00f75=176, 2\"),\n    @(\"132\u00f75=26, 2\", \"624\u00f72=312, 0\"),\n    @(\"623\u00f78=77, 7\", \"595\u00f72=297, 1\"),\n    @(\"528\u00f79=58, 6\", \"779\u00f73=259, 2\"),\n    @(\"465\u00f75=93, 0\", \"889\u00f77=127, 0\"),\n    @(\"335\u00f78=41, 7\", \"979\u00f73=326, 1\"),\n    @(\"313\u00f76=52, 1\", \"855\u00f79=95, 0\"),\n    @(\"956\u00f77=136, 4\", \"104\u00f77=14, 6\"),\n    @(\"442\u00f74=110, 2\", \"388\u00f72=194, 0\"),\n    @(\"227\u00f79=25, 2\", \"391\u00f72=195, 1\"),\n    @(\"471\u00f76=78, 3\", \"395\u00f79=43, 8\"),\n    @(\"267\u00f79=29, 6\", \"946\u00f72=473, 0\"),\n    @(\"858\u00f79=95, 3\", \"825\u00f72=412, 1\"),\n    @(\"583\u00f73=194, 1\", \"690\u00f79=76, 6\"),\n    @(\"564\u00f72=282, 0\", \"995\u00f77=142, 1\"),\n    @(\"824\u00f76=137, 2\", \"855\u00f79=95, 0\"),\n    @(\"390\u00f76=65, 0\", \"503\u00f76=83, 5\"),\n    @(\"290\u00f79=32, 2\", \"202\u00f72=101, 0\"),\n    @(\"580\u00f72=290, 0\", \"546\u00f78=68, 2\"),\n    @(\"556\u00f73=185, 1\", \"329\u00f79=36, 5\"),\n    @(\"303\u00f74=75, 3\", \"591\u00f78=73, 7\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
